$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: create rows 44 and 45 by copying formatting from existing rows ---
# This reuses existing style indices (no new cellXfs/border entries are created),
# matching the target workbook where styles.xml is untouched.
# Row 44
$ws.Range("A43").Copy($ws.Range("A44"))   # style 6 (bordered)
$ws.Range("B27").Copy($ws.Range("B44"))   # no explicit style
$ws.Range("C43").Copy($ws.Range("C44"))   # style 6
$ws.Range("D43").Copy($ws.Range("D44"))   # style 6
$ws.Range("E43").Copy($ws.Range("E44"))   # style 6
$ws.Range("F27").Copy($ws.Range("F44"))   # no explicit style

# Row 45
$ws.Range("A43").Copy($ws.Range("A45"))   # style 6
$ws.Range("B27").Copy($ws.Range("B45"))   # no explicit style
$ws.Range("C43").Copy($ws.Range("C45"))   # style 6
$ws.Range("D43").Copy($ws.Range("D45"))   # style 6
$ws.Range("E43").Copy($ws.Range("E45"))   # style 6
$ws.Range("F43").Copy($ws.Range("F45"))   # style 6

# --- Step 2: fill in the actual values, in the exact order the new unique
# strings were added to the shared-strings table upstream ---
$ws.Range("F44").Value = "gIT99L3TQ9s"
$ws.Range("E44").Value = "Northern Ireland"
$ws.Range("D44").Value = "Larne"
$ws.Range("B44").Value = "54.85621973585378, -5.815691135936153"
$ws.Range("A44").Value = "LIVE, CITY"

$ws.Range("F45").Value = "VhVgZi2lGv0"
$ws.Range("D45").Value = "Tallinn"
$ws.Range("E45").Value = "Estonia"
$ws.Range("C45").Value = "24/7 Live Stream 4K - Digital"
$ws.Range("C44").Value = "Larne, Northern Ireland"
$ws.Range("B45").Value = "59.437358306894886, 24.75017667351407"
$ws.Range("A45").Value = "LIVE, CITY, TRAFFIC"

# --- Step 3: update the saved selection to match the target workbook ---
$ws.Range("A46").Select()
